# Weekly update: a new week of Acelga price data (2021-09-10, serial 44449)
# is inserted right after the header block of this sub-series, pushing the
# existing rows 268-296 down to 271-299.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 fresh rows at the top of the data block (rows 268-270), shifting
# the existing rows 268-296 down to 271-299.
$ws.Rows("268:270").Insert()

# Row 268: Acelga, Extra
$ws.Range("A268").Value = 9
$ws.Range("B268").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C268").Value = "Metropolitana"
$ws.Range("D268").Value = 44449
$ws.Range("E268").Value = 13
$ws.Range("F268").Value = 100112009
$ws.Range("G268").Value = "Acelga"
$ws.Range("H268").Value = "Sin especificar"
$ws.Range("I268").Value = "Extra"
$ws.Range("J268").Value = 25
$ws.Range("K268").Value = 12000
$ws.Range("L268").Value = 13000
$ws.Range("M268").Value = 12520
$ws.Range("N268").Value = "$/docena de atados"
$ws.Range("O268").Value = "Región Metropolitana"
$ws.Range("P268").Value = 4173
$ws.Range("Q268").Value = 3
$ws.Range("R268").Value = "Hortaliza"

# Row 269: Acelga, Primera
$ws.Range("A269").Value = 9
$ws.Range("B269").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C269").Value = "Metropolitana"
$ws.Range("D269").Value = 44449
$ws.Range("E269").Value = 13
$ws.Range("F269").Value = 100112009
$ws.Range("G269").Value = "Acelga"
$ws.Range("H269").Value = "Sin especificar"
$ws.Range("I269").Value = "Primera"
$ws.Range("J269").Value = 52
$ws.Range("K269").Value = 10000
$ws.Range("L269").Value = 11000
$ws.Range("M269").Value = 10500
$ws.Range("N269").Value = "$/docena de atados"
$ws.Range("O269").Value = "Región Metropolitana"
$ws.Range("P269").Value = 3500
$ws.Range("Q269").Value = 3
$ws.Range("R269").Value = "Hortaliza"

# Row 270: Acelga, Segunda
$ws.Range("A270").Value = 9
$ws.Range("B270").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C270").Value = "Metropolitana"
$ws.Range("D270").Value = 44449
$ws.Range("E270").Value = 13
$ws.Range("F270").Value = 100112009
$ws.Range("G270").Value = "Acelga"
$ws.Range("H270").Value = "Sin especificar"
$ws.Range("I270").Value = "Segunda"
$ws.Range("J270").Value = 34
$ws.Range("K270").Value = 8000
$ws.Range("L270").Value = 9000
$ws.Range("M270").Value = 8500
$ws.Range("N270").Value = "$/docena de atados"
$ws.Range("O270").Value = "Región Metropolitana"
$ws.Range("P270").Value = 2833
$ws.Range("Q270").Value = 3
$ws.Range("R270").Value = "Hortaliza"
